$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.02492466666666666
$ws.Range("H2").Value = 0.07477399999999999
$ws.Range("I2").Value = 0.5549276228148621
$ws.Range("J2").Value = 0.6515968803102261
$ws.Range("M2").Value = 0.2687865
$ws.Range("N2").Value = 0.537573
$ws.Range("O2").Value = 0.01336253673722166
$ws.Range("P2").Value = 0.008961926404151501
$ws.Range("Q2").Value = 0.006699413916999999
$ws.Range("R2").Value = 0.040196483502
$ws.Range("S2").Value = 0.007415240746362681
$ws.Range("T2").Value = 0.005839563286514961
# Row 3
$ws.Range("G3").Value = 0.02492466666666666
$ws.Range("H3").Value = 0.07477399999999999
$ws.Range("I3").Value = 0.5549276228148621
$ws.Range("J3").Value = 0.6515968803102261
$ws.Range("O3").Value = 0.06358094663744968
$ws.Range("P3").Value = 0.06396327759577582
$ws.Range("Q3").Value = 0.03187681254955555
$ws.Range("S3").Value = 0.03528282357383854
$ws.Range("T3").Value = 0.0416782721358245
# Row 4
$ws.Range("G4").Value = 0.02492466666666666
$ws.Range("H4").Value = 0.07477399999999999
$ws.Range("I4").Value = 0.5549276228148621
$ws.Range("J4").Value = 0.6515968803102261
$ws.Range("M4").Value = 7.880035333333335
$ws.Range("N4").Value = 23.640106
$ws.Range("O4").Value = 0.3917505590209012
$ws.Range("P4").Value = 0.394106270512731
$ws.Range("Q4").Value = 0.1964072540048889
$ws.Range("R4").Value = 1.767665286044
$ws.Range("S4").Value = 0.217393206453862
$ws.Range("T4").Value = 0.2567984163767936
# Row 5
$ws.Range("G5").Value = 0.02492466666666666
$ws.Range("H5").Value = 0.07477399999999999
$ws.Range("I5").Value = 0.5549276228148621
$ws.Range("J5").Value = 0.6515968803102261
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = 0.0919155
$ws.Range("N5").Value = 0.183831
$ws.Range("O5").Value = 0.004569516123280365
$ws.Range("P5").Value = 0.003064662646378398
$ws.Range("Q5").Value = 0.002290963199
$ws.Range("R5").Value = 0.013745779194
$ws.Range("S5").Value = 0.002535750719706157
$ws.Range("T5").Value = 0.001996924619583446
# Row 6
$ws.Range("G6").Value = 0.02492466666666666
$ws.Range("H6").Value = 0.07477399999999999
$ws.Range("I6").Value = 0.5549276228148621
$ws.Range("J6").Value = 0.6515968803102261
$ws.Range("M6").Value = 2.558308333333333
$ws.Range("N6").Value = 7.674925
$ws.Range("O6").Value = 0.1271845464311154
$ws.Range("P6").Value = 0.1279493445678679
$ws.Range("Q6").Value = 0.06376498243888888
$ws.Range("R6").Value = 0.57388484195
$ws.Range("S6").Value = 0.07057821800980531
$ws.Range("T6").Value = 0.08337139375816091
# Row 7
$ws.Range("G7").Value = 0.02492466666666666
$ws.Range("H7").Value = 0.07477399999999999
$ws.Range("I7").Value = 0.5549276228148621
$ws.Range("J7").Value = 0.6515968803102261
$ws.Range("M7").Value = 8.036958666666665
$ws.Range("N7").Value = 24.110876
$ws.Range("O7").Value = 0.3995518950500317
$ws.Range("P7").Value = 0.4019545182730954
$ws.Range("Q7").Value = 0.2003185157804444
$ws.Range("R7").Value = 1.802866642024
$ws.Range("S7").Value = 0.2217223833112874
$ws.Range("T7").Value = 0.2619123101333487
# Row 8
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.0199905
$ws.Range("H8").Value = 0.039981
$ws.Range("I8").Value = 0.4450723771851379
$ws.Range("J8").Value = 0.3484031196897739
$ws.Range("M8").Value = 0.2687865
$ws.Range("N8").Value = 0.537573
$ws.Range("O8").Value = 0.01336253673722166
$ws.Range("P8").Value = 0.008961926404151501
$ws.Range("Q8").Value = 0.00537317652825
$ws.Range("R8").Value = 0.021492706113
$ws.Range("S8").Value = 0.005947295990858983
$ws.Range("T8").Value = 0.00312236311763654
# Row 9
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.0199905
$ws.Range("H9").Value = 0.039981
$ws.Range("I9").Value = 0.4450723771851379
$ws.Range("J9").Value = 0.3484031196897739
$ws.Range("O9").Value = 0.06358094663744968
$ws.Range("P9").Value = 0.06396327759577582
$ws.Range("Q9").Value = 0.0255663768665
$ws.Range("R9").Value = 0.153398261199
$ws.Range("S9").Value = 0.02829812306361113
$ws.Range("T9").Value = 0.02228500545995131
# Row 10
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.0199905
$ws.Range("H10").Value = 0.039981
$ws.Range("I10").Value = 0.4450723771851379
$ws.Range("J10").Value = 0.3484031196897739
$ws.Range("M10").Value = 7.880035333333335
$ws.Range("N10").Value = 23.640106
$ws.Range("O10").Value = 0.3917505590209012
$ws.Range("P10").Value = 0.394106270512731
$ws.Range("Q10").Value = 0.157525846331
$ws.Range("R10").Value = 0.9451550779860002
$ws.Range("S10").Value = 0.1743573525670392
$ws.Range("T10").Value = 0.1373078541359374
# Row 11
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.0199905
$ws.Range("H11").Value = 0.039981
$ws.Range("I11").Value = 0.4450723771851379
$ws.Range("J11").Value = 0.3484031196897739
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.5
$ws.Range("M11").Value = 0.0919155
$ws.Range("N11").Value = 0.183831
$ws.Range("O11").Value = 0.004569516123280365
$ws.Range("P11").Value = 0.003064662646378398
$ws.Range("Q11").Value = 0.00183743680275
$ws.Range("R11").Value = 0.007349747211000001
$ws.Range("S11").Value = 0.002033765403574208
$ws.Range("T11").Value = 0.001067738026794952
# Row 12
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.0199905
$ws.Range("H12").Value = 0.039981
$ws.Range("I12").Value = 0.4450723771851379
$ws.Range("J12").Value = 0.3484031196897739
$ws.Range("M12").Value = 2.558308333333333
$ws.Range("N12").Value = 7.674925
$ws.Range("O12").Value = 0.1271845464311154
$ws.Range("P12").Value = 0.1279493445678679
$ws.Range("Q12").Value = 0.0511418627375
$ws.Range("R12").Value = 0.306851176425
$ws.Range("S12").Value = 0.05660632842131007
$ws.Range("T12").Value = 0.04457795080970701
# Row 13
$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.0199905
$ws.Range("H13").Value = 0.039981
$ws.Range("I13").Value = 0.4450723771851379
$ws.Range("J13").Value = 0.3484031196897739
$ws.Range("M13").Value = 8.036958666666665
$ws.Range("N13").Value = 24.110876
$ws.Range("O13").Value = 0.3995518950500317
$ws.Range("P13").Value = 0.4019545182730954
$ws.Range("Q13").Value = 0.160662822226
$ws.Range("R13").Value = 0.963976933356
$ws.Range("S13").Value = 0.1778295117387444
$ws.Range("T13").Value = 0.1400422081397467
